$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = 53
$ws.Range("B28").Value = "11:26 added"
$ws.Range("C28").Value = "riya-morankar"
$ws.Range("D28").Value = "N/A"
$ws.Range("E28").Value = "edit1 to main"

$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "2025-06-20"
